$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the waiting time values in column C (result data)
$ws.Range("C2").Value = 30
$ws.Range("C3").Value = 20
$ws.Range("C4").Value = 24
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 15
$ws.Range("C7").Value = 20

# Update the selected cell/range to a single cell G2 (no more horizontal
# selection line spanning C2:C7)
$ws.Range("G2").Select()
